$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC)
$ws.Range("H33").Value = 7548.5
$ws.Range("I33").Value = 126.71429
$ws.Range("J33").Value = 14970.286
$ws.Range("K33").Value = 126.71429
$ws.Range("L33").Value = 14970.286
$ws.Range("M33").Value = 102.28571
$ws.Range("N33").Value = -15428.286

# Row 76 (ALC)
$ws.Range("H76").Value = 3733.1333
$ws.Range("I76").Value = 3644.1667
$ws.Range("J76").Value = 4089
$ws.Range("K76").Value = 3644.1667
$ws.Range("L76").Value = 4089
$ws.Range("M76").Value = -3329.1667
$ws.Range("N76").Value = -4719

# Row 79 (ALC)
$ws.Range("H79").Value = 3733.1333
$ws.Range("I79").Value = 3644.1667
$ws.Range("J79").Value = 4089
$ws.Range("K79").Value = 3644.1667
$ws.Range("L79").Value = 4089
$ws.Range("M79").Value = -2552.1667
$ws.Range("N79").Value = -6273

# Row 132 (ALC)
$ws.Range("H132").Value = 3125.2273
$ws.Range("I132").Value = 3010.4707
$ws.Range("J132").Value = 3515.4
$ws.Range("K132").Value = 9031.4121
$ws.Range("L132").Value = 10546.2
$ws.Range("M132").Value = -6501.4121
$ws.Range("N132").Value = -15606.2

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (BSM)
$ws.Range("H82").Value = 16194.077
$ws.Range("J82").Value = 26074.715
$ws.Range("L82").Value = 26074.715
$ws.Range("N82").Value = -26840.715

# Row 85 (BSM)
$ws.Range("H85").Value = 16194.077
$ws.Range("J85").Value = 26074.715
$ws.Range("L85").Value = 26074.715
$ws.Range("N85").Value = -28726.715

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 1712.8572
$ws.Range("I16").Value = 1331.6666
$ws.Range("K16").Value = 1331.6666
$ws.Range("M16").Value = -1044.6666

# Row 31 (CRP)
$ws.Range("H31").Value = 2470.4443
$ws.Range("I31").Value = 2039.3914
$ws.Range("J31").Value = 3233.077
$ws.Range("K31").Value = 2039.3914
$ws.Range("L31").Value = 3233.077
$ws.Range("M31").Value = -1744.3914
$ws.Range("N31").Value = -3823.077

# Row 34 (CRP)
$ws.Range("H34").Value = 2470.4443
$ws.Range("I34").Value = 2039.3914
$ws.Range("J34").Value = 3233.077
$ws.Range("K34").Value = 2039.3914
$ws.Range("L34").Value = 3233.077
$ws.Range("M34").Value = -1837.3914
$ws.Range("N34").Value = -3637.077

# Row 99 (CRP)
$ws.Range("H99").Value = 1854.8695
$ws.Range("I99").Value = 1907.3055
$ws.Range("J99").Value = 1666.1
$ws.Range("K99").Value = 1907.3055
$ws.Range("L99").Value = 1666.1
$ws.Range("M99").Value = -409.3054999999999
$ws.Range("N99").Value = -4662.1

# Row 113 (CRP)
$ws.Range("H113").Value = 1712.8572
$ws.Range("I113").Value = 1331.6666
$ws.Range("K113").Value = 1331.6666
$ws.Range("M113").Value = 838.3334

# Row 126 (CRP)
$ws.Range("H126").Value = 1854.8695
$ws.Range("I126").Value = 1907.3055
$ws.Range("J126").Value = 1666.1
$ws.Range("K126").Value = 5721.916499999999
$ws.Range("L126").Value = 4998.299999999999
$ws.Range("M126").Value = -3251.916499999999
$ws.Range("N126").Value = -9938.299999999999

# Row 132 (CRP)
$ws.Range("H132").Value = 2261.4285
$ws.Range("I132").Value = 1551.3334
$ws.Range("J132").Value = 3539.6
$ws.Range("K132").Value = 4654.0002
$ws.Range("L132").Value = 10618.8
$ws.Range("M132").Value = -2124.0002
$ws.Range("N132").Value = -15678.8

$ws = $wb.Worksheets.Item("CUL")
# Row 107 (CUL)
$ws.Range("H107").Value = 1143.3334
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1143.3334
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3430.0002
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7270.0002

# Row 113 (CUL)
$ws.Range("H113").Value = 767.2069
$ws.Range("I113").Value = 474.4
$ws.Range("J113").Value = 885.2742
$ws.Range("K113").Value = 1423.2
$ws.Range("L113").Value = 2655.8226
$ws.Range("M113").Value = 746.8000000000002
$ws.Range("N113").Value = -6995.8226

# Row 131 (CUL)
$ws.Range("H131").Value = 5000844
$ws.Range("I131").Value = 1051.1111
$ws.Range("J131").Value = 5495329
$ws.Range("K131").Value = 3153.3333
$ws.Range("L131").Value = 16485987
$ws.Range("M131").Value = 1886.6667
$ws.Range("N131").Value = -16496067

# Row 132 (CUL)
$ws.Range("H132").Value = 1164.7059
$ws.Range("I132").Value = 1084.1666
$ws.Range("J132").Value = 1255.3125
$ws.Range("K132").Value = 9757.499400000001
$ws.Range("L132").Value = 11297.8125
$ws.Range("M132").Value = -7227.499400000001
$ws.Range("N132").Value = -16357.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 1699557
$ws.Range("I70").Value = 2987977
$ws.Range("J70").Value = 8505.875
$ws.Range("K70").Value = 2987977
$ws.Range("L70").Value = 8505.875
$ws.Range("M70").Value = -2987707
$ws.Range("N70").Value = -9045.875

# Row 73 (GSM)
$ws.Range("H73").Value = 1699557
$ws.Range("I73").Value = 2987977
$ws.Range("J73").Value = 8505.875
$ws.Range("K73").Value = 2987977
$ws.Range("L73").Value = 8505.875
$ws.Range("M73").Value = -2987041
$ws.Range("N73").Value = -10377.875

# Row 80 (GSM)
$ws.Range("H80").Value = 4090
$ws.Range("I80").Value = 1750
$ws.Range("J80").Value = 4675
$ws.Range("K80").Value = 1750
$ws.Range("L80").Value = 4675
$ws.Range("M80").Value = -752
$ws.Range("N80").Value = -6671

# Row 83 (GSM)
$ws.Range("H83").Value = 4090
$ws.Range("I83").Value = 1750
$ws.Range("J83").Value = 4675
$ws.Range("K83").Value = 8750
$ws.Range("L83").Value = 23375
$ws.Range("M83").Value = -3758
$ws.Range("N83").Value = -33359

# Row 113 (GSM)
$ws.Range("H113").Value = 111112536
$ws.Range("I113").Value = 1364
$ws.Range("J113").Value = 166668130
$ws.Range("K113").Value = 1364
$ws.Range("L113").Value = 166668130
$ws.Range("M113").Value = 806
$ws.Range("N113").Value = -166672470

# Row 126 (GSM)
$ws.Range("H126").Value = 20835854
$ws.Range("I126").Value = 4022.6667
$ws.Range("J126").Value = 33334952
$ws.Range("K126").Value = 12068.0001
$ws.Range("L126").Value = 100004856
$ws.Range("M126").Value = -9598.000100000001
$ws.Range("N126").Value = -100009796

# Row 132 (GSM)
$ws.Range("H132").Value = 7352.857
$ws.Range("I132").Value = 7895.3887
$ws.Range("K132").Value = 23686.1661
$ws.Range("M132").Value = -21156.1661

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 3471.4285
$ws.Range("I61").Value = 1700
$ws.Range("J61").Value = 5833.3335
$ws.Range("K61").Value = 1700
$ws.Range("L61").Value = 5833.3335
$ws.Range("M61").Value = -1498
$ws.Range("N61").Value = -6237.3335

# Row 113 (LTW)
$ws.Range("H113").Value = 3471.4285
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 5833.3335
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 5833.3335
$ws.Range("M113").Value = 470
$ws.Range("N113").Value = -10173.3335

# Row 132 (LTW)
$ws.Range("H132").Value = 5632.7085
$ws.Range("I132").Value = 5947.0527
$ws.Range("J132").Value = 4438.2
$ws.Range("K132").Value = 17841.1581
$ws.Range("L132").Value = 13314.6
$ws.Range("M132").Value = -15311.1581
$ws.Range("N132").Value = -18374.6

# Row 136 (LTW)
$ws.Range("H136").Value = 3198.25
$ws.Range("I136").Value = 3798
$ws.Range("K136").Value = 11394
$ws.Range("M136").Value = -8844

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 27031176
$ws.Range("I132").Value = 33334574
$ws.Range("J132").Value = 16614.285
$ws.Range("K132").Value = 100003722
$ws.Range("L132").Value = 49842.855
$ws.Range("M132").Value = -100001192
$ws.Range("N132").Value = -54902.855

# Row 135 (WVR)
$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

# Row 136 (WVR)
$ws.Range("H136").Value = 2125
$ws.Range("I136").Value = 2650
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 7950
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -5400
$ws.Range("N136").Value = -9900

# Row 137 (WVR)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
